# Add a new row of user data (row 6) to the Users sheet, matching the
# existing email/username/password columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "test3@gmail.com"
$ws.Range("B6").Value = "Subhadra Mahato"

# Store the password as text (not a number) so it matches the existing
# "inline string" representation used elsewhere in the sheet, then clear
# the temporary number-format override so no extra style is left behind.
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "12345678"
$ws.Range("C6").ClearFormats()
